{"js": "// The document contains a 2-digit x 2-digit multiplication practice table.\n// Each cell's equation text is replaced with its updated equation, matched\n// by exact old-text -> new-text pairs (same order the values appear in the body).\nconst replacements = [\n  [\"89\u00d758=5162\", \"90\u00d775=6750\"],\n  [\"84\u00d741=3444\", \"14\u00d771=994\"],\n  [\"54\u00d755=2970\", \"21\u00d768=1428\"],\n  [\"18\u00d777=1386\", \"92\u00d778=7176\"],\n  [\"31\u00d721=651\", \"14\u00d798=1372\"],\n  [\"90\u00d742=3780\", \"73\u00d786=6278\"],\n  [\"83\u00d759=4897\", \"92\u00d733=3036\"],\n  [\"15\u00d749=735\", \"31\u00d752=1612\"],\n  [\"25\u00d790=2250\", \"43\u00d770=3010\"],\n  [\"24\u00d746=1104\", \"33\u00d745=1485\"],\n  [\"83\u00d725=2075\", \"91\u00d743=3913\"],\n  [\"87\u00d731=2697\", \"17\u00d711=187\"],\n  [\"11\u00d754=594\", \"54\u00d761=3294\"],\n  [\"47\u00d795=4465\", \"57\u00d785=4845\"],\n  [\"16\u00d799=1584\", \"74\u00d761=4514\"],\n  [\"59\u00d735=2065\", \"57\u00d749=2793\"],\n  [\"24\u00d793=2232\", \"38\u00d786=3268\"],\n  [\"17\u00d735=595\", \"42\u00d742=1764\"],\n  [\"43\u00d767=2881\", \"35\u00d724=840\"],\n  [\"33\u00d758=1914\", \"41\u00d736=1476\"],\n  [\"40\u00d722=880\", \"39\u00d767=2613\"],\n  [\"46\u00d720=920\", \"93\u00d776=7068\"],\n  [\"58\u00d798=5684\", \"99\u00d797=9603\"],\n  [\"51\u00d780=4080\", \"89\u00d772=6408\"],\n  [\"22\u00d751=1122\", \"69\u00d783=5727\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`No match found for \"${oldText}\"`);\n  }\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update each two-digit multiplication equation in the practice table\n# to its new value, matching old text to new text one-to-one in document order.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"89\u00d758=5162\"; New = \"90\u00d775=6750\" }\n    @{ Old = \"84\u00d741=3444\"; New = \"14\u00d771=994\" }\n    @{ Old = \"54\u00d755=2970\"; New = \"21\u00d768=1428\" }\n    @{ Old = \"18\u00d777=1386\"; New = \"92\u00d778=7176\" }\n    @{ Old = \"31\u00d721=651\"; New = \"14\u00d798=1372\" }\n    @{ Old = \"90\u00d742=3780\"; New = \"73\u00d786=6278\" }\n    @{ Old = \"83\u00d759=4897\"; New = \"92\u00d733=3036\" }\n    @{ Old = \"15\u00d749=735\"; New = \"31\u00d752=1612\" }\n    @{ Old = \"25\u00d790=2250\"; New = \"43\u00d770=3010\" }\n    @{ Old = \"24\u00d746=1104\"; New = \"33\u00d745=1485\" }\n    @{ Old = \"83\u00d725=2075\"; New = \"91\u00d743=3913\" }\n    @{ Old = \"87\u00d731=2697\"; New = \"17\u00d711=187\" }\n    @{ Old = \"11\u00d754=594\"; New = \"54\u00d761=3294\" }\n    @{ Old = \"47\u00d795=4465\"; New = \"57\u00d785=4845\" }\n    @{ Old = \"16\u00d799=1584\"; New = \"74\u00d761=4514\" }\n    @{ Old = \"59\u00d735=2065\"; New = \"57\u00d749=2793\" }\n    @{ Old = \"24\u00d793=2232\"; New = \"38\u00d786=3268\" }\n    @{ Old = \"17\u00d735=595\"; New = \"42\u00d742=1764\" }\n    @{ Old = \"43\u00d767=2881\"; New = \"35\u00d724=840\" }\n    @{ Old = \"33\u00d758=1914\"; New = \"41\u00d736=1476\" }\n    @{ Old = \"40\u00d722=880\"; New = \"39\u00d767=2613\" }\n    @{ Old = \"46\u00d720=920\"; New = \"93\u00d776=7068\" }\n    @{ Old = \"58\u00d798=5684\"; New = \"99\u00d797=9603\" }\n    @{ Old = \"51\u00d780=4080\"; New = \"89\u00d772=6408\" }\n    @{ Old = \"22\u00d751=1122\"; New = \"69\u00d783=5727\" }\n)\n\nforeach ($pair in $replacements) {\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    # Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards, MatchSoundsLike,\n    #         MatchAllWordForms, Forward, Wrap, Format, ReplaceWith, Replace)\n    # Wrap=1 -> wdFindContinue, Replace=2 -> wdReplaceAll\n    $find.Execute($pair.Old, $false, $true, $false, $false, $false, $true, 1, $false, $pair.New, 2) | Out-Null\n}\n"}
